$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 53 (Serie index 55) with revised figures
$ws.Range("B53").Value = 65221884
$ws.Range("D53").Value = 54423

# Update existing row 54 (Serie index 56) with revised figures
$ws.Range("B54").Value = 67667925
$ws.Range("D54").Value = 51271

# Append new row 55 for period 01-04-2021
# Use a leading apostrophe so Excel stores the literal text instead of
# auto-converting the date-like string to a serial date, then reset the
# cell style back to Normal (matches the other "Serie" text cells, which
# carry no explicit style).
$ws.Range("A55").Value = "'01-04-2021"
$ws.Range("A55").Style = "Normal"
$ws.Range("B55").Value = 71982736
$ws.Range("C55").Value = 71933370
$ws.Range("D55").Value = 49367
